$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: corrected Hydrogen "Iron & steel" value, clear the stray
# "Non-metallic minerals" value that no longer applies.
$ws.Range("B3").Value = 51373247.57370247
$ws.Range("D3").ClearContents()

# Row 7 was "Other" -> now it holds the "Biogas" figures, with an
# updated value.
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 29022.15988383117

# New row 8 takes over as "Other", carrying the same row styling as
# the rest of column A's labels.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 21616.80632522705
